$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.8958887032678717
$ws.Cells.Item(2, 3).Value = 0.2189747176543904
$ws.Cells.Item(2, 4).Value = 0.01005867474103894
$ws.Cells.Item(2, 6).Value = 0.5415983124126313
$ws.Cells.Item(2, 7).Value = 0.392050031227015
$ws.Cells.Item(2, 8).Value = 0.5071163709027644
$ws.Cells.Item(2, 14).Value = 0.8462759456202136
$ws.Cells.Item(3, 2).Value = 0.7881322745201373
$ws.Cells.Item(3, 3).Value = 0.1926122297886366
$ws.Cells.Item(3, 4).Value = 0.009208595059675417
$ws.Cells.Item(3, 6).Value = 0.5304946671969333
$ws.Cells.Item(3, 7).Value = 0.3813551720773347
$ws.Cells.Item(3, 8).Value = 0.5073623226635959
$ws.Cells.Item(3, 14).Value = 0.8537561001142393
$ws.Cells.Item(4, 2).Value = 0.7219650120557617
$ws.Cells.Item(4, 3).Value = 0.176349297489395
$ws.Cells.Item(4, 4).Value = 0.008684786428574398
$ws.Cells.Item(4, 6).Value = 0.5242024040731579
$ws.Cells.Item(4, 7).Value = 0.3752425414654113
$ws.Cells.Item(4, 8).Value = 0.5079074215363022
$ws.Cells.Item(4, 14).Value = 0.8587893591180205
$ws.Cells.Item(5, 2).Value = 0.6950009142082934
$ws.Cells.Item(5, 3).Value = 0.1697030218368809
$ws.Cells.Item(5, 4).Value = 0.008470880999460917
$ws.Cells.Item(5, 6).Value = 0.5217697533767449
$ws.Cells.Item(5, 7).Value = 0.3728649808215891
$ws.Cells.Item(5, 8).Value = 0.5082283899329099
$ws.Cells.Item(5, 14).Value = 0.8609512304883609
$ws.Cells.Item(6, 2).Value = 0.6905235429206016
$ws.Cells.Item(6, 3).Value = 0.1685982706208335
$ws.Cells.Item(6, 4).Value = 0.008435335559454415
$ws.Cells.Item(6, 6).Value = 0.5213737374856251
$ws.Cells.Item(6, 7).Value = 0.3724770138014719
$ws.Cells.Item(6, 8).Value = 0.508287648381696
$ws.Cells.Item(6, 14).Value = 0.86131689932472
$ws.Cells.Item(7, 2).Value = 0.7216013654925746
$ws.Cells.Item(7, 3).Value = 0.1762597401782386
$ws.Cells.Item(7, 4).Value = 0.008681903419859083
$ws.Cells.Item(7, 6).Value = 0.5241690648933925
$ws.Cells.Item(7, 7).Value = 0.3752100188027754
$ws.Cells.Item(7, 8).Value = 0.5079113503900601
$ws.Cells.Item(7, 14).Value = 0.858818066294603
$ws.Cells.Item(8, 2).Value = 0.8587355807603672
$ws.Cells.Item(8, 3).Value = 0.2099008572948549
$ws.Cells.Item(8, 4).Value = 0.009765962784726412
$ws.Cells.Item(8, 6).Value = 0.5376603443206065
$ws.Cells.Item(8, 7).Value = 0.3882677246916302
$ws.Cells.Item(8, 8).Value = 0.5071192251972718
$ws.Cells.Item(8, 14).Value = 0.8487637120507543
$ws.Cells.Item(9, 2).Value = 1.127603479416109
$ws.Cells.Item(9, 3).Value = 0.2752612790830256
$ws.Cells.Item(9, 4).Value = 0.01187640516798893
$ws.Cells.Item(9, 6).Value = 0.5683160143899499
$ws.Cells.Item(9, 7).Value = 0.4175140647109856
$ws.Cells.Item(9, 8).Value = 0.5087052291547565
$ws.Cells.Item(9, 14).Value = 0.8325407286362179
$ws.Cells.Item(10, 2).Value = 1.325106220210671
$ws.Cells.Item(10, 3).Value = 0.3229083116866605
$ws.Cells.Item(10, 4).Value = 0.01341682562393487
$ws.Cells.Item(10, 6).Value = 0.5934439358486827
$ws.Cells.Item(10, 7).Value = 0.4412750650270283
$ws.Cells.Item(10, 8).Value = 0.5118027172632083
$ws.Cells.Item(10, 14).Value = 0.822750969856969
$ws.Cells.Item(11, 2).Value = 1.414948575068991
$ws.Cells.Item(11, 3).Value = 0.344503104032583
$ws.Cells.Item(11, 4).Value = 0.01411526172640976
$ws.Cells.Item(11, 6).Value = 0.6054508907480312
$ws.Cells.Item(11, 7).Value = 0.4525901635999645
$ws.Cells.Item(11, 8).Value = 0.5136354665463045
$ws.Cells.Item(11, 14).Value = 0.8187597660038435
$ws.Cells.Item(12, 2).Value = 1.448968767213444
$ws.Cells.Item(12, 3).Value = 0.3526688668275995
$ws.Cells.Item(12, 4).Value = 0.01437939428215174
$ws.Cells.Item(12, 6).Value = 0.6100811955836605
$ws.Cells.Item(12, 7).Value = 0.4569485780243951
$ws.Cells.Item(12, 8).Value = 0.5143907166752371
$ws.Cells.Item(12, 14).Value = 0.8173148821372678
$ws.Cells.Item(13, 2).Value = 1.441641973963897
$ws.Cells.Item(13, 3).Value = 0.3509107479738134
$ws.Cells.Item(13, 4).Value = 0.01432252446298321
$ws.Cells.Item(13, 6).Value = 0.609080250606155
$ws.Cells.Item(13, 7).Value = 0.4560066270575334
$ws.Cells.Item(13, 8).Value = 0.5142253321029528
$ws.Cells.Item(13, 14).Value = 0.8176231055904495
$ws.Cells.Item(14, 2).Value = 1.417747464394608
$ws.Cells.Item(14, 3).Value = 0.3451751422099392
$ws.Cells.Item(14, 4).Value = 0.01413699918867906
$ws.Cells.Item(14, 6).Value = 0.6058301506975994
$ws.Cells.Item(14, 7).Value = 0.4529472523772995
$ws.Cells.Item(14, 8).Value = 0.5136963725990853
$ws.Cells.Item(14, 14).Value = 0.8186395615957807
$ws.Cells.Item(15, 2).Value = 1.403111213554155
$ws.Cells.Item(15, 3).Value = 0.3416603863203136
$ws.Cells.Item(15, 4).Value = 0.01402331348440811
$ws.Cells.Item(15, 6).Value = 0.6038502692407803
$ws.Cells.Item(15, 7).Value = 0.4510829115286015
$ws.Cells.Item(15, 8).Value = 0.5133803522756892
$ws.Cells.Item(15, 14).Value = 0.8192708314131991
$ws.Cells.Item(16, 2).Value = 1.319234691287591
$ws.Cells.Item(16, 3).Value = 0.3214954125701865
$ws.Cells.Item(16, 4).Value = 0.01337113317502769
$ws.Cells.Item(16, 6).Value = 0.592670904902036
$ws.Cells.Item(16, 7).Value = 0.4405458538833216
$ws.Cells.Item(16, 8).Value = 0.5116914942474011
$ws.Cells.Item(16, 14).Value = 0.823021107924383
$ws.Cells.Item(17, 2).Value = 1.267777930756665
$ws.Cells.Item(17, 3).Value = 0.3091042000672815
$ws.Cells.Item(17, 4).Value = 0.01297043742088277
$ws.Cells.Item(17, 6).Value = 0.5859607759267647
$ws.Cells.Item(17, 7).Value = 0.4342119433719631
$ws.Cells.Item(17, 8).Value = 0.5107641654818167
$ws.Cells.Item(17, 14).Value = 0.8254401852060838
$ws.Cells.Item(18, 2).Value = 1.238181161443663
$ws.Cells.Item(18, 3).Value = 0.3019695765924553
$ws.Cells.Item(18, 4).Value = 0.012739751559252
$ws.Cells.Item(18, 6).Value = 0.5821554750858411
$ws.Cells.Item(18, 7).Value = 0.4306164309636102
$ws.Cells.Item(18, 8).Value = 0.5102706517359934
$ws.Cells.Item(18, 14).Value = 0.8268750771071254
$ws.Cells.Item(19, 2).Value = 1.228160193121766
$ws.Cells.Item(19, 3).Value = 0.2995526292565671
$ws.Cells.Item(19, 4).Value = 0.01266160870079602
$ws.Cells.Item(19, 6).Value = 0.5808763516011481
$ws.Cells.Item(19, 7).Value = 0.429407199017902
$ws.Cells.Item(19, 8).Value = 0.5101103934805167
$ws.Cells.Item(19, 14).Value = 0.8273683775521974
$ws.Cells.Item(20, 2).Value = 1.273255617296343
$ws.Cells.Item(20, 3).Value = 0.3104240461433676
$ws.Cells.Item(20, 4).Value = 0.01301311468133548
$ws.Cells.Item(20, 6).Value = 0.58666946741792
$ws.Cells.Item(20, 7).Value = 0.4348812678532141
$ws.Cells.Item(20, 8).Value = 0.5108587535352882
$ws.Cells.Item(20, 14).Value = 0.8251781677132612
$ws.Cells.Item(21, 2).Value = 1.424765898664646
$ws.Cells.Item(21, 3).Value = 0.3468601482265115
$ws.Cells.Item(21, 4).Value = 0.01419150210082876
$ws.Cells.Item(21, 6).Value = 0.6067825113487402
$ws.Cells.Item(21, 7).Value = 0.4538438597049179
$ws.Cells.Item(21, 8).Value = 0.513850076710952
$ws.Cells.Item(21, 14).Value = 0.8183391987799453
$ws.Cells.Item(22, 2).Value = 1.523779520737151
$ws.Cells.Item(22, 3).Value = 0.3706048632056422
$ws.Cells.Item(22, 4).Value = 0.01495959791590451
$ws.Cells.Item(22, 6).Value = 0.6204147520569734
$ws.Cells.Item(22, 7).Value = 0.4666665916373063
$ws.Cells.Item(22, 8).Value = 0.5161621024807914
$ws.Cells.Item(22, 14).Value = 0.8142571610082356
$ws.Cells.Item(23, 2).Value = 1.470934993021956
$ws.Cells.Item(23, 3).Value = 0.3579381849517915
$ws.Cells.Item(23, 4).Value = 0.01454984405425819
$ws.Cells.Item(23, 6).Value = 0.6130941633912244
$ws.Cells.Item(23, 7).Value = 0.4597832687028642
$ws.Cells.Item(23, 8).Value = 0.5148953630719149
$ws.Cells.Item(23, 14).Value = 0.8164003410169443
$ws.Cells.Item(24, 2).Value = 1.270779197495813
$ws.Cells.Item(24, 3).Value = 0.3098273772090181
$ws.Cells.Item(24, 4).Value = 0.01299382128872395
$ws.Cells.Item(24, 6).Value = 0.586348904693196
$ws.Cells.Item(24, 7).Value = 0.4345785232297175
$ws.Cells.Item(24, 8).Value = 0.510815866880904
$ws.Cells.Item(24, 14).Value = 0.8252964883492524
$ws.Cells.Item(25, 2).Value = 1.054873467308596
$ws.Cells.Item(25, 3).Value = 0.25764499605026
$ws.Cells.Item(25, 4).Value = 0.01130719971978777
$ws.Cells.Item(25, 6).Value = 0.5595686096682755
$ws.Cells.Item(25, 7).Value = 0.4092068087597056
$ws.Cells.Item(25, 8).Value = 0.5079381837395118
$ws.Cells.Item(25, 14).Value = 0.8365555421716451
